$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.141.54'
$ws.Range('E2').Value = '  -1.95%  '
$ws.Range('D3').Value = '1.852.90'
$ws.Range('E3').Value = '  -0.88%  '
$ws.Range('E4').Value = '  +0.06%  '
$origStyle = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '239.54'
$ws.Range('D5').Style = $origStyle
$ws.Range('E5').Value = '  -0.60%  '
$origStyle = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6926'
$ws.Range('D6').Style = $origStyle
$ws.Range('E6').Value = '  -4.89%  '
$ws.Range('E7').Value = '  +0.05%  '
$origStyle = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07726'
$ws.Range('D8').Style = $origStyle
$ws.Range('E8').Value = '  +9.15%  '
$origStyle = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3044'
$ws.Range('D9').Style = $origStyle
$origStyle = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.48'
$ws.Range('D10').Style = $origStyle
$ws.Range('E10').Value = '  -3.64%  '
$origStyle = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08142'
$ws.Range('D11').Style = $origStyle
$ws.Range('E11').Value = '  -1.19%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.861.78'
$ws.Range('E12').Value = '  -2.67%  '
$ws.Range('B13').Value = 'Polygon'
$ws.Range('C13').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$origStyle = $ws.Range('D13').Style
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.7273'
$ws.Range('D13').Style = $origStyle
$ws.Range('E13').Value = '  -2.45%  '
$origStyle = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.237'
$ws.Range('D14').Style = $origStyle
$ws.Range('E14').Value = '  -1.54%  '
$origStyle = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '89.17'
$ws.Range('D15').Style = $origStyle
$ws.Range('E15').Value = '  -3.29%  '
$ws.Range('D16').Value = '29.158.46'
$ws.Range('E16').Value = '  -2.01%  '
$origStyle = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.804'
$ws.Range('D17').Style = $origStyle
$ws.Range('E17').Value = '  -3.60%  '
$origStyle = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000007784'
$ws.Range('D18').Style = $origStyle
$ws.Range('E18').Value = '  -0.20%  '
$origStyle = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.19'
$ws.Range('D19').Style = $origStyle
$ws.Range('E19').Value = '  -1.17%  '
$origStyle = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '237.04'
$ws.Range('D20').Style = $origStyle
$ws.Range('E20').Value = '  -4.19%  '
$origStyle = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.000'
$ws.Range('D21').Style = $origStyle
$ws.Range('E21').Value = '  +0.10%  '
$ws.Range('D22').Value = '2.098.51'
$ws.Range('E22').Value = '  -1.56%  '
$ws.Range('E23').Value = '  +0.05%  '
$origStyle = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.639'
$ws.Range('D24').Style = $origStyle
$ws.Range('E24').Value = '  -0.87%  '
$origStyle = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.022'
$ws.Range('D25').Style = $origStyle
$ws.Range('E25').Value = '  -1.45%  '
$origStyle = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '161.46'
$ws.Range('D26').Style = $origStyle
$ws.Range('E26').Value = '  -0.88%  '
$origStyle = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1449'
$ws.Range('D27').Style = $origStyle
$ws.Range('E27').Value = '  -5.45%  '
$origStyle = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.10'
$ws.Range('D28').Style = $origStyle
$ws.Range('E28').Value = '  -2.31%  '
$origStyle = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.983'
$ws.Range('D29').Style = $origStyle
$ws.Range('E29').Value = '  -1.79%  '
$ws.Range('E30').Value = '  -2.28%  '
$origStyle = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.483'
$ws.Range('D31').Style = $origStyle
$ws.Range('E31').Value = '  -1.16%  '
$origStyle = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.497'
$ws.Range('D32').Style = $origStyle
$ws.Range('E32').Value = '  -2.03%  '
$origStyle = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.017'
$ws.Range('D33').Style = $origStyle
$ws.Range('E33').Value = '  -4.10%  '
$ws.Range('E34').Value = '  -0.93%  '
$ws.Range('E35').Value = '  -3.48%  '
$origStyle = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.038'
$ws.Range('D36').Style = $origStyle
$ws.Range('E36').Value = '  +3.78%  '
$origStyle = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.7045'
$ws.Range('D37').Style = $origStyle
$ws.Range('E37').Value = '  -6.52%  '
$origStyle = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.658'
$ws.Range('D38').Style = $origStyle
$ws.Range('E38').Value = '  -1.62%  '
$origStyle = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01862'
$ws.Range('D39').Style = $origStyle
$ws.Range('E39').Value = '  -3.49%  '
$origStyle = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.677'
$ws.Range('D40').Style = $origStyle
$ws.Range('E40').Value = '  -2.46%  '
$origStyle = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9375'
$ws.Range('D41').Style = $origStyle
$ws.Range('E41').Value = '  +8.07%  '
$origStyle = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.048'
$ws.Range('D42').Style = $origStyle
$ws.Range('E42').Value = '  +0.91%  '
$ws.Range('D43').Value = '1.080.80'
$ws.Range('E43').Value = '  +1.51%  '
$origStyle = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.4285'
$ws.Range('D44').Style = $origStyle
$ws.Range('E44').Value = '  -4.14%  '
$origStyle = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '70.47'
$ws.Range('D45').Style = $origStyle
$ws.Range('E45').Value = '  -1.19%  '
$origStyle = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.000'
$ws.Range('D46').Style = $origStyle
$ws.Range('E46').Value = '  -0.01%  '
$origStyle = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '103.28'
$ws.Range('D47').Style = $origStyle
$ws.Range('E47').Value = '  -0.74%  '
$origStyle = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.786'
$ws.Range('D48').Style = $origStyle
$ws.Range('E48').Value = '  -1.96%  '
$ws.Range('D49').Value = '1.996.65'
$ws.Range('E49').Value = '  -1.58%  '
$origStyle = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.225'
$ws.Range('D50').Style = $origStyle
$ws.Range('E50').Value = '  -3.23%  '
$origStyle = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.028'
$ws.Range('D51').Style = $origStyle
$ws.Range('E51').Value = '  -6.07%  '
